$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format before writing, so numeric-looking strings
# (e.g. "299.04") are preserved exactly as text instead of being parsed as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '42.132.90'
$ws.Range('E2').Value = '  -1.03%  '

$ws.Range('D3').Value = '2.272.88'
$ws.Range('E3').Value = '  -1.35%  '

$ws.Range('E4').Value = '  -0.09%  '

$ws.Range('D5').Value = '299.04'
$ws.Range('E5').Value = '  -1.19%  '

$ws.Range('D6').Value = '95.53'
$ws.Range('E6').Value = '  -3.91%  '

$ws.Range('E7').Value = '  -2.55%  '

$ws.Range('E8').Value = '  -0.04%  '

$ws.Range('D9').Value = '0.490'
$ws.Range('E9').Value = '  -2.55%  '

$ws.Range('D10').Value = '33.05'
$ws.Range('E10').Value = '  -4.80%  '

$ws.Range('D11').Value = '0.0787'
$ws.Range('E11').Value = '  -0.67%  '

$ws.Range('D12').Value = '48.15'
$ws.Range('E12').Value = '  -6.76%  '

$ws.Range('D14').Value = '15.96'
$ws.Range('E14').Value = '  +2.01%  '

$ws.Range('D15').Value = '6.67'
$ws.Range('E15').Value = '  -0.92%  '

$ws.Range('D16').Value = '2.624.97'
$ws.Range('E16').Value = '  -1.43%  '

$ws.Range('D17').Value = '2.277.32'
$ws.Range('E17').Value = '  -2.35%  '

$ws.Range('D18').Value = '0.786'
$ws.Range('E18').Value = '  -2.32%  '

$ws.Range('D19').Value = '42.085.12'
$ws.Range('E19').Value = '  -1.03%  '

$ws.Range('D20').Value = '11.70'
$ws.Range('E20').Value = '  +2.14%  '

$ws.Range('D21').Value = '0.0₃0890'
$ws.Range('E21').Value = '  -1.73%  '

$ws.Range('D22').Value = '5.98'
$ws.Range('E22').Value = '  -1.56%  '

$ws.Range('D23').Value = '66.27'
$ws.Range('E23').Value = '  -2.27%  '

$ws.Range('D24').Value = '235.30'
$ws.Range('E24').Value = '  +0.19%  '

$ws.Range('E25').Value = '  -0.73%  '

$ws.Range('E26').Value = '  +0.14%  '

$ws.Range('E27').Value = '  -2.31%  '

$ws.Range('D28').Value = '23.78'
$ws.Range('E28').Value = '  -4.63%  '

$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').Value = '2.17'
$ws.Range('E29').Value = '  -5.47%  '

$ws.Range('B30').Value = 'Monero'
$ws.Range('C30').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D30').Value = '167.85'
$ws.Range('E30').Value = '  +2.41%  '

$ws.Range('B31').Value = 'InjectiveProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D31').Value = '33.67'
$ws.Range('E31').Value = '  -2.93%  '

$ws.Range('B32').Value = 'Cosmos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D32').Value = '9.14'
$ws.Range('E32').Value = '  -0.02%  '

$ws.Range('E33').Value = '  -0.07%  '

$ws.Range('D34').Value = '4.70'
$ws.Range('E34').Value = '  +6.11%  '

$ws.Range('D35').Value = '4.90'
$ws.Range('E35').Value = '  -2.11%  '

$ws.Range('B36').Value = 'Celestia'
$ws.Range('C36').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D36').Value = '16.76'
$ws.Range('E36').Value = '  +0.21%  '

$ws.Range('B37').Value = 'WEMIXToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D37').Value = '2.35'
$ws.Range('E37').Value = '  -2.84%  '

$ws.Range('D38').Value = '0.0687'
$ws.Range('E38').Value = '  -2.32%  '

$ws.Range('D39').Value = '2.80'
$ws.Range('E39').Value = '  -2.78%  '

$ws.Range('D40').Value = '0.0984'
$ws.Range('E40').Value = '  -1.68%  '

$ws.Range('D41').Value = '0.109'
$ws.Range('E41').Value = '  -2.17%  '

$ws.Range('D42').Value = '1.72'
$ws.Range('E42').Value = '  -4.40%  '

$ws.Range('D43').Value = '2.30'
$ws.Range('E43').Value = '  -5.80%  '

$ws.Range('D44').Value = '1.958.16'
$ws.Range('E44').Value = '  -0.37%  '

$ws.Range('E45').Value = '  -1.39%  '

$ws.Range('D46').Value = '17.62'
$ws.Range('E46').Value = '  -4.66%  '

$ws.Range('D47').Value = '9.55'
$ws.Range('E47').Value = '  -6.48%  '

$ws.Range('D48').Value = '2.77'
$ws.Range('E48').Value = '  -4.08%  '

$ws.Range('D49').Value = '2.494.96'
$ws.Range('E49').Value = '  -1.34%  '

$ws.Range('B50').Value = 'MultiversX'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D50').Value = '51.95'
$ws.Range('E50').Value = '  -6.63%  '

$ws.Range('B51').Value = 'HuobiToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D51').Value = '2.74'
$ws.Range('E51').Value = '  -3.67%  '

# Restore the default (Normal) style on column D so no stray number format
# is left behind on cells, matching the original workbook formatting.
$ws.Range("D2:D51").Style = "Normal"
